# Insert a new weekly record at row 253 (this pushes the existing rows
# 253-351 down to 254-352, preserving all their data/formatting), then
# populate the newly inserted row with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(253).Insert()

$ws.Range("A253").Value = 6
$ws.Range("B253").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C253").Value = "Metropolitana"
$ws.Range("D253").Value = 44468
$ws.Range("E253").Value = 13
$ws.Range("F253").Value = 100112044
$ws.Range("G253").Value = "Perejil"
$ws.Range("H253").Value = "Sin especificar"
$ws.Range("I253").Value = "Primera"
$ws.Range("J253").Value = 210
$ws.Range("K253").Value = 7000
$ws.Range("L253").Value = 8000
$ws.Range("M253").Value = 7381
$ws.Range("N253").Value = "`$/docena de atados"
$ws.Range("O253").Value = "Región Metropolitana"
$ws.Range("P253").Value = 2460
$ws.Range("Q253").Value = 3
$ws.Range("R253").Value = "Hortaliza"
